$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-fill R3:R48 with SUM(Bn:Qn) as a single range formula assignment so Excel
# stores it as one shared-formula group (t="shared" si="0"), matching the
# "fill down" edit captured in the diff.
$ws.Range("R3:R48").Formula = "=SUM(B3:Q3)"

# Row 49 totals switch from SUM(...) to COUNT(...). B49 was edited on its own
# (stays a plain formula) while C49:Q49 were filled together as one shared
# group (t="shared" si="1"), and R49 is a brand new total cell.
$ws.Range("B49").Formula = "=COUNT(B3:B48)"
$ws.Range("C49:Q49").Formula = "=COUNT(C3:C48)"
$ws.Range("R49").Formula = "=SUM(R3:R48)"

# Selection moved to the whole of row 39.
$ws.Rows(39).Select()

# Window was minimized/resized in the session that produced this save.
try {
    $excel.ActiveWindow.WindowState = -4140
    $excel.ActiveWindow.Left = 5385
    $excel.ActiveWindow.Top = 4095
    $excel.ActiveWindow.Width = 7500
    $excel.ActiveWindow.Height = 6000
} catch {}
